# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, shifting the existing Late / heading /
# Outstanding columns one place to the right, then leave the
# "Repayment schedule" tab selected/active (matching the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture the width of column M (the column immediately to the left of the
# insertion point) so the freshly inserted column inherits the same width,
# the way Excel's own "Insert Column" behaves.
$leftWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new blank column at N, pushing Late/heading/Outstanding to O/P/Q.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet and restore the last selected
# cell on it.
$ws.Select()
$ws.Range("R8").Select()
